# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" sheet right after "总计" (and right before the
# existing "2022-Q3" sheet), shifting "2022-Q3" and "2021-Q1" one slot to
# the right, and updates the summary ("总计") sheet to include the new
# quarter's row.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q3Sheet    = $wb.Worksheets.Item("2022-Q3")

# xlPasteFormats - used below to copy formatting only, without disturbing
# the values/types already written into a cell.
$xlPasteFormats = -4122

# A scratch area (far outside the used range of any sheet) that always has
# the plain, un-styled default format; used as a formatting source so text
# values typed with a leading apostrophe don't keep the "number stored as
# text" quote-prefix style.
$blankFormatCell = $totalSheet.Range("Z100")

# --- 1. Create the new "2022-Q4" sheet -------------------------------------
# Copy the "2022-Q3" sheet (placing the copy right before it, i.e. right
# after "总计") so the new sheet inherits the same header/row styling, then
# overwrite its data with the Q4 figures.
$q3Sheet.Copy($q3Sheet, $null)
$q4Sheet = $wb.Worksheets.Item("2022-Q3 (2)")
$q4Sheet.Name = "2022-Q4"

# Data rows. Numeric-looking codes/percentages are stored as TEXT (matching
# the other quarterly sheets): write them with a leading apostrophe so
# Excel doesn't auto-coerce them to numbers, then strip the resulting
# quote-prefix formatting so the cell ends up plain (un-styled) text.
$q4Sheet.Range("A2").Value = 0
$q4Sheet.Range("B2").Value = "'513690"
$q4Sheet.Range("C2").Value = "博时恒生港股通高股息率ETF"
$q4Sheet.Range("D2").Value = "'5.20"
$q4Sheet.Range("E2").Value = "'98.04"
$q4Sheet.Range("F2").Value = "'2.67"
$q4Sheet.Range("G2").Value = "'0.1388"
$q4Sheet.Range("H2").Value = 6

$q4Sheet.Range("A3").Value = 1
$q4Sheet.Range("B3").Value = "'159726"
$q4Sheet.Range("C3").Value = "华夏恒生中国内地企业高股息率ETF"
$q4Sheet.Range("D3").Value = "'0.84"
$q4Sheet.Range("E3").Value = "'98.34"
$q4Sheet.Range("F3").Value = "'2.91"
$q4Sheet.Range("G3").Value = "'0.0244"
$q4Sheet.Range("H3").Value = 5

$q4Sheet.Range("A4").Value = 2
$q4Sheet.Range("B4").Value = "'005702"
$q4Sheet.Range("C4").Value = "恒生前海港股通高股息低波动指数"
$q4Sheet.Range("D4").Value = "'0.23"
$q4Sheet.Range("E4").Value = "'94.47"
$q4Sheet.Range("F4").Value = "'2.53"
$q4Sheet.Range("G4").Value = "'0.0058"
$q4Sheet.Range("H4").Value = 4

# Column A keeps the bordered/centered header style used by row 2 (copied
# from "2022-Q3"); propagate it to the newly-added rows 3-4.
$q4Sheet.Range("A2").Copy()
$q4Sheet.Range("A3:A4").PasteSpecial($xlPasteFormats)

# Strip the quote-prefix formatting picked up from the leading apostrophes.
$blankFormatCell.Copy()
$q4Sheet.Range("B2:B4").PasteSpecial($xlPasteFormats)
$q4Sheet.Range("D2:G4").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# --- 2. Update the "总计" (summary) sheet -----------------------------------
# Old row 2 ("2022-Q3"/1/0) becomes the new "2022-Q4"/3/0.17 row, and the
# previous rows 2-3 shift down to make room for the original "2022-Q3" row.
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.17

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q1"
$totalSheet.Range("C4").Value = 2
$totalSheet.Range("D4").Value = 0.03

# Column A keeps the bordered/centered style; propagate it to new row 4.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# Restore the originally-active tab ("2021-Q1") — copying/renaming sheets
# above shifts Excel's selection to whichever sheet was touched last.
$wb.Worksheets.Item("2021-Q1").Activate()
